$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 306
$ws.Cells.Item(2, 2).Value = 28
$ws.Cells.Item(2, 3).Value = "楊*晟"
$ws.Cells.Item(2, 4).Value = "2023-06-28 10:10:00"
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = "IN"
